$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "A4" = -20.575
    "C4" = -12.566
    "A6" = -22.275
    "D6" = -7.810999999999998
    "A7" = -19.948
    "D7" = -8.537000000000001
    "A8" = -21.903
    "C8" = -12.453
    "D8" = -8.094000000000001
    "C9" = -11
    "D10" = -7.641999999999999
    "C12" = -11.011
    "D13" = -7.817
    "D14" = -8.306000000000001
    "A16" = -22.025
    "D16" = -8.086
    "C17" = -12.948
    "C18" = -12.448
    "C19" = -12.126
    "A20" = -19.828
    "C20" = -11.729
    "A21" = -20.269
    "C26" = -12.5
    "A28" = -22.113
    "A29" = -21.324
    "A30" = -21.758
    "D30" = -7.206999999999999
    "C31" = -12.711
    "A32" = -21.837
    "D37" = -8.440999999999999
    "C39" = -12.117
    "A40" = -19.927
    "C40" = -12.038
    "D40" = -7.996
    "C41" = -12.117
    "C42" = -12.14
    "C43" = -12.485
    "D44" = -7.766
    "A46" = -21.89
    "C47" = -12.125
    "C48" = -12.359
    "A51" = -21.71
    "A52" = -21.972
    "C54" = -13.268
    "A57" = -22.192
    "A59" = -22.42
    "A62" = -21.873
    "C62" = -12.286
    "C63" = -11.046
    "C64" = -10.674
    "A66" = -21.627
    "D70" = -7.702
    "A73" = -20.318
    "A74" = -21.036
    "C76" = -12.45
    "A77" = -20.407
    "C81" = -12.187
    "C84" = -13.162
    "C89" = -12.901
    "D89" = -8.123000000000001
    "D91" = -7.580000000000001
    "A92" = -21.722
    "D93" = -7.465999999999999
    "C94" = -11.87
    "D98" = -7.008
    "A100" = -22.159
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
